$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the styling of the existing header row (reuse H1's exact format)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows 2-20 for columns I (I0) and J (IF)
$data = @(
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(9, 9),
    @(9, 9),
    @(8, 8),
    @(8, 8),
    @(7, 8),
    @(10, 10),
    @(7, 7),
    @(7, 8),
    @(8, 9),
    @(3, 4),
    @(8, 9),
    @(1, 1),
    @(7, 7),
    @(5, 6),
    @(8, 8),
    @(7, 7)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
